{"js": "// 1) \"...cubren por completo y de manera precisa la necesidad.\" ->\n//    \"...cubren por completo y de manera precisa las necesidades.\"\nconst needNecesidad = context.document.body.search(\"precisa la necesidad.\", { matchCase: true });\nneedNecesidad.load(\"text\");\nawait context.sync();\n\nif (needNecesidad.items.length > 0) {\n  needNecesidad.items[0].insertText(\n    \"precisa las necesidades.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 2) \"...dentro del sistema.\" ->\n//    \"...dentro del sistema, y genere autom\u00e1ticamente documentos editables.\"\nconst needSistema = context.document.body.search(\"dentro del sistema.\", { matchCase: true });\nneedSistema.load(\"text\");\nawait context.sync();\n\nif (needSistema.items.length > 0) {\n  needSistema.items[0].insertText(\n    \"dentro del sistema, y genere autom\u00e1ticamente documentos editables.\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"...cubren por completo y de manera precisa la necesidad.\" ->\n#    \"...cubren por completo y de manera precisa las necesidades.\"\n$rng1 = $d.Content\n$rng1.Find.MatchCase = $true\n$rng1.Find.MatchWholeWord = $false\n$found1 = $rng1.Find.Execute(\"precisa la necesidad.\")\nif ($found1) {\n    $rng1.Text = \"precisa las necesidades.\"\n}\n\n# 2) \"...dentro del sistema.\" ->\n#    \"...dentro del sistema, y genere autom\u00e1ticamente documentos editables.\"\n$rng2 = $d.Content\n$rng2.Find.MatchCase = $true\n$rng2.Find.MatchWholeWord = $false\n$found2 = $rng2.Find.Execute(\"dentro del sistema.\")\nif ($found2) {\n    $rng2.Text = \"dentro del sistema, y genere autom\u00e1ticamente documentos editables.\"\n}\n"}
